# Regenerate the "K" (strikeouts) column values in the save_data sheet.
# The source data was re-scraped using the actual strikeout count (K)
# instead of the previous "Strike#" (total pitches that were strikes),
# so column G (header "K") needs its values replaced for every game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 2
    4  = 1
    5  = 1
    6  = 8
    7  = 7
    8  = 8
    9  = 3
    10 = 7
    11 = 6
    12 = 6
    13 = 9
    14 = 6
    15 = 5
    16 = 6
    17 = 6
    18 = 4
    19 = 4
    20 = 5
    21 = 4
    22 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
